$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("M64").ClearContents()
$ws.Range("H64").Value = 7580
$ws.Range("K64").Value = 0
$ws.Range("I64").Value = 0
# Row 67
$ws.Range("M67").ClearContents()
$ws.Range("H67").Value = 7580
$ws.Range("K67").Value = 0
$ws.Range("I67").Value = 0
# Row 74
$ws.Range("H74").Value = 12771.286
$ws.Range("I74").Value = 6079.8
$ws.Range("K74").Value = 6079.8
$ws.Range("M74").Value = -5143.8
# Row 77
$ws.Range("K77").Value = 30399
$ws.Range("M77").Value = -25719
$ws.Range("I77").Value = 6079.8
$ws.Range("H77").Value = 12771.286
# Row 132
$ws.Range("M132").Value = -32241.422
$ws.Range("H132").Value = 11950.75
$ws.Range("J132").Value = 13319.8
$ws.Range("N132").Value = -45019.39999999999
$ws.Range("L132").Value = 39959.39999999999
$ws.Range("I132").Value = 11590.474
$ws.Range("K132").Value = 34771.422
# Row 135
$ws.Range("M135").Value = -6078.8181
$ws.Range("K135").Value = 8613.8181
$ws.Range("I135").Value = 957.0909
$ws.Range("H135").Value = 1127.3334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H37").Value = 99999
# Row 45
$ws.Range("I45").Value = 2372.0908
$ws.Range("K45").Value = 2372.0908
$ws.Range("H45").Value = 2861.3125
$ws.Range("M45").Value = -1995.0908
# Row 122
$ws.Range("H122").Value = 2249.5
$ws.Range("M122").Value = -2050
$ws.Range("K122").Value = 4500
$ws.Range("I122").Value = 1500
# Row 135
$ws.Range("J135").Value = 79999
$ws.Range("L135").Value = 79999
$ws.Range("N135").Value = -90139
$ws.Range("H135").Value = 79999
# Row 137
$ws.Range("H137").Value = 46806
$ws.Range("L137").Value = 75000
$ws.Range("M137").Value = -27609
$ws.Range("N137").Value = -85200
$ws.Range("I137").Value = 32709
$ws.Range("K137").Value = 32709
$ws.Range("J137").Value = 75000

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("K36").Value = 1143.8334
$ws.Range("M36").Value = -609.8334
$ws.Range("I36").Value = 1143.8334
$ws.Range("H36").Value = 1266.1428
# Row 86
$ws.Range("K86").Value = 2770.8
$ws.Range("M86").Value = -1647.8
$ws.Range("H86").Value = 6285.4
$ws.Range("I86").Value = 2770.8
# Row 89
$ws.Range("I89").Value = 2770.8
$ws.Range("H89").Value = 6285.4
$ws.Range("M89").Value = -8238
$ws.Range("K89").Value = 13854
# Row 105
$ws.Range("K105").Value = 9243993
$ws.Range("I105").Value = 9243993
$ws.Range("M105").Value = -9242246
$ws.Range("H105").Value = 5839074

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("N58").Value = -5603.25
$ws.Range("K58").Value = 1893.7142
$ws.Range("J58").Value = 5197.25
$ws.Range("I58").Value = 1893.7142
$ws.Range("M58").Value = -1690.7142
$ws.Range("L58").Value = 5197.25
# Row 132
$ws.Range("M132").Value = -3741.600199999999
$ws.Range("H132").Value = 2090.5334
$ws.Range("I132").Value = 2090.5334
$ws.Range("K132").Value = 6271.600199999999
# Row 134
$ws.Range("K134").Value = 7404.249899999999
$ws.Range("L134").Value = 0
$ws.Range("H134").Value = 2468.0833
$ws.Range("I134").Value = 2468.0833
$ws.Range("J134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("M134").Value = -4869.249899999999
# Row 136
$ws.Range("M136").Value = -3131.142599999999
$ws.Range("L136").Value = 15591.75
$ws.Range("J136").Value = 5197.25
$ws.Range("K136").Value = 5681.142599999999
$ws.Range("I136").Value = 1893.7142
$ws.Range("N136").Value = -20691.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("I2").Value = 16.75
$ws.Range("H2").Value = 20.363636
$ws.Range("K2").Value = 100.5
$ws.Range("M2").Value = 12.5
# Row 35
$ws.Range("N35").Value = -1761
$ws.Range("L35").Value = 1185
$ws.Range("J35").Value = 395
$ws.Range("H35").Value = 395
# Row 59
$ws.Range("K59").Value = 2700
$ws.Range("I59").Value = 900
$ws.Range("L59").Value = 2250
$ws.Range("M59").Value = -2160
$ws.Range("N59").Value = -3330
$ws.Range("J59").Value = 750
$ws.Range("H59").Value = 787.5
# Row 75
$ws.Range("N75").Value = -23428
$ws.Range("L75").Value = 21432
$ws.Range("J75").Value = 7144
$ws.Range("H75").Value = 5781.8
# Row 78
$ws.Range("J78").Value = 7144
$ws.Range("N78").Value = -74280
$ws.Range("L78").Value = 64296
$ws.Range("H78").Value = 5781.8
# Row 80
$ws.Range("J80").Value = 5130.8
$ws.Range("N80").Value = -17264.4
$ws.Range("L80").Value = 15392.4
$ws.Range("H80").Value = 4297.9565
# Row 83
$ws.Range("H83").Value = 4297.9565
$ws.Range("J83").Value = 5130.8
$ws.Range("N83").Value = -55537.2
$ws.Range("L83").Value = 46177.2
# Row 88
$ws.Range("L88").Value = 37500
$ws.Range("N88").Value = -38356
$ws.Range("H88").Value = 12500
$ws.Range("J88").Value = 12500
# Row 91
$ws.Range("L91").Value = 37500
$ws.Range("H91").Value = 12500
$ws.Range("J91").Value = 12500
$ws.Range("N91").Value = -40464
# Row 92
$ws.Range("J92").Value = 2048
$ws.Range("L92").Value = 6144
$ws.Range("N92").Value = -8640
$ws.Range("H92").Value = 1548.8462
# Row 97
$ws.Range("H97").Value = 215.4
$ws.Range("N97").Value = -1484.75
$ws.Range("J97").Value = 164.25
$ws.Range("L97").Value = 492.75
# Row 117
$ws.Range("N117").Value = -8884.499900000001
$ws.Range("H117").Value = 592.2857
$ws.Range("L117").Value = 2000.4999
$ws.Range("J117").Value = 666.8333

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("N2").Value = -741.75
$ws.Range("L2").Value = 515.75
$ws.Range("J2").Value = 515.75
$ws.Range("H2").Value = 222.56522
# Row 80
$ws.Range("J80").Value = 4999
$ws.Range("N80").Value = -6995
$ws.Range("L80").Value = 4999
$ws.Range("H80").Value = 4999
# Row 83
$ws.Range("H83").Value = 4999
$ws.Range("J83").Value = 4999
$ws.Range("N83").Value = -34979
$ws.Range("L83").Value = 24995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("I7").Value = 1167.6666
$ws.Range("L7").Value = 0
$ws.Range("H7").Value = 1167.6666
$ws.Range("K7").Value = 1167.6666
$ws.Range("N7").ClearContents()
$ws.Range("J7").Value = 0
$ws.Range("M7").Value = -1055.6666
# Row 30
$ws.Range("M30").Value = -497.44446
$ws.Range("K30").Value = 605.44446
$ws.Range("I30").Value = 605.44446
$ws.Range("H30").Value = 605.44446
# Row 93
$ws.Range("L93").Value = 4749
$ws.Range("K93").Value = 4049.2
$ws.Range("I93").Value = 4049.2
$ws.Range("N93").Value = -7245
$ws.Range("J93").Value = 4749
$ws.Range("H93").Value = 4165.8335
$ws.Range("M93").Value = -2801.2
# Row 126
$ws.Range("K126").Value = 3502.9998
$ws.Range("L126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("H126").Value = 1167.6666
$ws.Range("I126").Value = 1167.6666
$ws.Range("N126").ClearContents()
$ws.Range("M126").Value = -1032.9998

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("K126").Value = 7576.3638
$ws.Range("L126").Value = 19831.5
$ws.Range("J126").Value = 6610.5
$ws.Range("H126").Value = 4245.4736
$ws.Range("I126").Value = 2525.4546
$ws.Range("N126").Value = -24771.5
$ws.Range("M126").Value = -5106.3638
# Row 132
$ws.Range("M132").Value = -6002.299999999999
$ws.Range("H132").Value = 2870.0833
$ws.Range("I132").Value = 2844.1
$ws.Range("K132").Value = 8532.299999999999
# Row 136
$ws.Range("M136").Value = -2214.706200000001
$ws.Range("H136").Value = 2383.32
$ws.Range("K136").Value = 4764.706200000001
$ws.Range("I136").Value = 1588.2354
